$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Cfh", "Itgam", "FAPs", 3, 1, 1.362097333333333, 4.086292, 0.0196292037450248, 0.0196292037450248, 1, 0.3333333333333333, 0.142723, 0.428169, 0.000470790275435748, 0.0004707902754357479, 0.1944026177053334, 1.749623559348, [double]"9.241238237704643E-06", [double]"9.241238237704638E-06")
for ($i = 0; $i -lt $row2.Length; $i++) { $ws.Cells.Item(2, $i + 1).Value = $row2[$i] }

$row3 = @("ECs", "Cfh", "Itgam", "M1", 3, 1, 1.362097333333333, 4.086292, 0.0196292037450248, 0.0196292037450248, 3, 1, 148.0881626666667, 444.264488, 0.4884879584272602, 0.4884879584272601, 201.7104914664996, 1815.394423198496, 0.009588629662959895, 0.009588629662959891)
for ($i = 0; $i -lt $row3.Length; $i++) { $ws.Cells.Item(3, $i + 1).Value = $row3[$i] }

$row4 = @("ECs", "Cfh", "Itgam", "M2", 3, 1, 1.362097333333333, 4.086292, 0.0196292037450248, 0.0196292037450248, 3, 1, 154.9253336666667, 464.776001, 0.5110412512973043, 0.5110412512973042, 211.0233838531436, 1899.210454678292, 0.0100313328438272, 0.0100313328438272)
for ($i = 0; $i -lt $row4.Length; $i++) { $ws.Cells.Item(4, $i + 1).Value = $row4[$i] }

$row5 = @("FAPs", "Cfh", "Itgam", "FAPs", 3, 1, 24.05951033333333, 72.17853099999999, 0.3467219403350491, 0.346721940335049, 1, 0.3333333333333333, 0.142723, 0.428169, 0.000470790275435748, 0.0004707902754357479, 3.433845493304334, 30.904609439739, 0.0001632333177899547, 0.0001632333177899547)
for ($i = 0; $i -lt $row5.Length; $i++) { $ws.Cells.Item(5, $i + 1).Value = $row5[$i] }

$row6 = @("FAPs", "Cfh", "Itgam", "M1", 3, 1, 24.05951033333333, 72.17853099999999, 0.3467219403350491, 0.346721940335049, 3, 1, 148.0881626666667, 444.264488, 0.4884879584272602, 0.4884879584272601, 3562.928679923014, 32066.35811930713, 0.1693694927762064, 0.1693694927762064)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, $i + 1).Value = $row6[$i] }

$row7 = @("FAPs", "Cfh", "Itgam", "M2", 3, 1, 24.05951033333333, 72.17853099999999, 0.3467219403350491, 0.346721940335049, 3, 1, 154.9253336666667, 464.776001, 0.5110412512973043, 0.5110412512973042, 3727.427666248281, 33546.84899623453, 0.1771892142410527, 0.1771892142410527)
for ($i = 0; $i -lt $row7.Length; $i++) { $ws.Cells.Item(7, $i + 1).Value = $row7[$i] }

$row8 = @("M1", "Cfh", "Itgam", "FAPs", 3, 1, 14.711625, 44.13487499999999, 0.2120094338917045, 0.2120094338917045, 1, 0.3333333333333333, 0.142723, 0.428169, 0.000470790275435748, 0.0004707902754357479, 2.099687254875, 18.897185293875, [double]"9.981197977685258E-05", [double]"9.981197977685256E-05")
for ($i = 0; $i -lt $row8.Length; $i++) { $ws.Cells.Item(8, $i + 1).Value = $row8[$i] }

$row9 = @("M1", "Cfh", "Itgam", "M1", 3, 1, 14.711625, 44.13487499999999, 0.2120094338917045, 0.2120094338917045, 3, 1, 148.0881626666667, 444.264488, 0.4884879584272602, 0.4884879584272601, 2178.617516091, 19607.557644819, 0.1035640555290779, 0.1035640555290779)
for ($i = 0; $i -lt $row9.Length; $i++) { $ws.Cells.Item(9, $i + 1).Value = $row9[$i] }

$row10 = @("M1", "Cfh", "Itgam", "M2", 3, 1, 14.711625, 44.13487499999999, 0.2120094338917045, 0.2120094338917045, 3, 1, 154.9253336666667, 464.776001, 0.5110412512973043, 0.5110412512973042, 2279.203411903874, 20512.83070713487, 0.1083455663828498, 0.1083455663828498)
for ($i = 0; $i -lt $row10.Length; $i++) { $ws.Cells.Item(10, $i + 1).Value = $row10[$i] }

$row11 = @("M2", "Cfh", "Itgam", "FAPs", 3, 1, 27.616616, 82.84984800000001, 0.3979834399099074, 0.3979834399099074, 1, 0.3333333333333333, 0.142723, 0.428169, 0.000470790275435748, 0.0004707902754357479, 3.941526285368001, 35.47373656831201, 0.0001873667332940518, 0.0001873667332940517)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, $i + 1).Value = $row11[$i] }

$row12 = @("M2", "Cfh", "Itgam", "M1", 3, 1, 27.616616, 82.84984800000001, 0.3979834399099074, 0.3979834399099074, 3, 1, 148.0881626666667, 444.264488, 0.4884879584272602, 0.4884879584272601, 4089.69392251087, 36807.24530259783, 0.1944101180494489, 0.1944101180494488)
for ($i = 0; $i -lt $row12.Length; $i++) { $ws.Cells.Item(12, $i + 1).Value = $row12[$i] }

$row13 = @("M2", "Cfh", "Itgam", "M2", 3, 1, 27.616616, 82.84984800000001, 0.3979834399099074, 0.3979834399099074, 3, 1, 154.9253336666667, 464.776001, 0.5110412512973043, 0.5110412512973042, 4278.513448544206, 38506.62103689785, 0.2033859551271646, 0.2033859551271645)
for ($i = 0; $i -lt $row13.Length; $i++) { $ws.Cells.Item(13, $i + 1).Value = $row13[$i] }

$row14 = @("sCs", "Cfh", "Itgam", "FAPs", 3, 1, 1.641521, 4.924563000000001, 0.02365598211831425, 0.02365598211831425, 1, 0.3333333333333333, 0.142723, 0.428169, 0.000470790275435748, 0.0004707902754357479, 0.2342828016830001, 2.108545215147001, [double]"1.113700633718429E-05", [double]"1.113700633718429E-05")
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, $i + 1).Value = $row14[$i] }

$row15 = @("sCs", "Cfh", "Itgam", "M1", 3, 1, 1.641521, 4.924563000000001, 0.02365598211831425, 0.02365598211831425, 3, 1, 148.0881626666667, 444.264488, 0.4884879584272602, 0.4884879584272601, 243.0898288687494, 2187.808459818745, 0.0115556624095671, 0.0115556624095671)
for ($i = 0; $i -lt $row15.Length; $i++) { $ws.Cells.Item(15, $i + 1).Value = $row15[$i] }

$row16 = @("sCs", "Cfh", "Itgam", "M2", 3, 1, 1.641521, 4.924563000000001, 0.02365598211831425, 0.02365598211831425, 3, 1, 154.9253336666667, 464.776001, 0.5110412512973043, 0.5110412512973042, 254.3131886458404, 2288.818697812563, 0.01208918270240997, 0.01208918270240996)
for ($i = 0; $i -lt $row16.Length; $i++) { $ws.Cells.Item(16, $i + 1).Value = $row16[$i] }
